$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) column values are stored as TEXT in this sheet even when they
# look numeric (e.g. "582.59"), and some use "." as a thousands separator
# (e.g. "65.754.09"). Force text so Excel does not auto-convert/reparse them,
# then restore the original (unstyled) cell style so no stray number-format
# style gets attached.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range('D2') '65.754.09'
$ws.Range('E2').Value = '  -3.28%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.478.44'
$ws.Range('E3').Value = '  -0.23%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
Set-TextValue $ws.Range('D5') '582.59'
$ws.Range('E5').Value = '  -2.02%  '

# Row 6
Set-TextValue $ws.Range('D6') '173.47'
$ws.Range('E6').Value = '  -4.94%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('E8').Value = '  -3.17%  '

# Row 9
Set-TextValue $ws.Range('D9') '3.478.83'
$ws.Range('E9').Value = '  -0.17%  '

# Row 10
$ws.Range('E10').Value = '  -7.05%  '

# Row 11
$ws.Range('E11').Value = '  -2.37%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.410'
$ws.Range('E12').Value = '  -4.89%  '

# Row 13
Set-TextValue $ws.Range('D13') '4.081.80'
$ws.Range('E13').Value = '  -0.01%  '

# Row 14
$ws.Range('E14').Value = '  +0.22%  '

# Row 15
Set-TextValue $ws.Range('D15') '30.04'
$ws.Range('E15').Value = '  -6.95%  '

# Row 16
Set-TextValue $ws.Range('D16') '65.937.75'
$ws.Range('E16').Value = '  -3.08%  '

# Row 17
$ws.Range('E17').Value = '  -3.61%  '

# Row 18
Set-TextValue $ws.Range('D18') '3.481.76'
$ws.Range('E18').Value = '  -0.12%  '

# Row 19
Set-TextValue $ws.Range('D19') '5.94'
$ws.Range('E19').Value = '  -4.54%  '

# Row 20
Set-TextValue $ws.Range('D20') '13.92'
$ws.Range('E20').Value = '  -1.48%  '

# Row 21
Set-TextValue $ws.Range('D21') '364.98'
$ws.Range('E21').Value = '  -7.76%  '

# Row 22
Set-TextValue $ws.Range('D22') '7.76'
$ws.Range('E22').Value = '  -2.30%  '

# Row 23
$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range('D23') '72.65'
$ws.Range('E23').Value = '  +0.48%  '

# Row 24
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D24') '1.00'
$ws.Range('E24').Value = '  +0.33%  '

# Row 25
Set-TextValue $ws.Range('D25') '0.534'
$ws.Range('E25').Value = '  -1.30%  '

# Row 26
Set-TextValue $ws.Range('D26') '0.0000123'
$ws.Range('E26').Value = '  +0.92%  '

# Row 27
Set-TextValue $ws.Range('D27') '9.69'
$ws.Range('E27').Value = '  -7.38%  '

# Row 28
$ws.Range('E28').Value = '  +0.15%  '

# Row 29
$ws.Range('E29').Value = '  -0.01%  '

# Row 30
Set-TextValue $ws.Range('D30') '24.06'
$ws.Range('E30').Value = '  +1.79%  '

# Row 31
Set-TextValue $ws.Range('D31') '5.78'
$ws.Range('E31').Value = '  -5.94%  '

# Row 32
$ws.Range('E32').Value = '  -4.07%  '

# Row 33
Set-TextValue $ws.Range('D33') '1.00'
$ws.Range('E33').Value = '  +0.04%  '

# Row 34
Set-TextValue $ws.Range('D34') '7.13'
$ws.Range('E34').Value = '  -3.39%  '

# Row 35
$ws.Range('E35').Value = '  -8.28%  '

# Row 36
Set-TextValue $ws.Range('D36') '1.54'
$ws.Range('E36').Value = '  -2.36%  '

# Row 37
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D37') '159.82'
$ws.Range('E37').Value = '  -1.28%  '

# Row 38
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D38') '29.40'
$ws.Range('E38').Value = '  +12.18%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.887'
$ws.Range('E39').Value = '  -0.70%  '

# Row 40
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D40') '2.816.07'
$ws.Range('E40').Value = '  +2.38%  '

# Row 41
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D41') '1.77'
$ws.Range('E41').Value = '  -5.77%  '

# Row 42
$ws.Range('E42').Value = '  -9.56%  '

# Row 43
Set-TextValue $ws.Range('D43') '4.46'
$ws.Range('E43').Value = '  -4.79%  '

# Row 44
Set-TextValue $ws.Range('D44') '6.36'
$ws.Range('E44').Value = '  -6.42%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.0686'
$ws.Range('E45').Value = '  -4.99%  '

# Row 46
Set-TextValue $ws.Range('D46') '39.92'
$ws.Range('E46').Value = '  -3.82%  '

# Row 47
Set-TextValue $ws.Range('D47') '24.15'
$ws.Range('E47').Value = '  -8.18%  '

# Row 48
$ws.Range('E48').Value = '  -3.92%  '

# Row 49
Set-TextValue $ws.Range('D49') '307.06'
$ws.Range('E49').Value = '  -6.98%  '

# Row 50
Set-TextValue $ws.Range('D50') '0.821'
$ws.Range('E50').Value = '  -3.45%  '

# Row 51
$ws.Range('E51').Value = '  -4.21%  '
